$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Hello"
$ws.Range("B22").Value = ","
$ws.Range("C22").Value = "World!!!"

$ws.Range("A23").Value = "ahmet"
$ws.Range("B23").Value = "can"
$ws.Range("C23").Value = "ozyurek"
